$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# --- Row 59: new "deadline" column D value, recalcs C/E automatically ---
$ws.Range("D59").Value = 21276

# --- Row 60: new commit row "tiny fixes" ---
$ws.Range("A60").Value = "tiny fixes"
$ws.Range("B60").Value = 20304
$ws.Range("D60").Value = 21276

# --- Row 61: new commit row "#61 track directive support" ---
$ws.Range("A61").Value = "#61 track directive support"
$ws.Range("B61").Value = 20312
$ws.Range("D61").Value = 21276
$ws.Range("F61").Value = "at least it was a tiny change"

# --- Rows 62-72: fill in the deadline/target column D with the same constant ---
$ws.Range("D62:D72").Value = 21276

# --- Restore the view/scroll state as closely as possible ---
$excel.ActiveWindow.ScrollRow = 42
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A61").Select()
